# "Beitrag Aargauer Kuratorium 2023"
#
# 1) Einnahmen sheet: append a new income row for the "Aargauer Kuratorium"
#    contribution to the Table3 listobject.
# 2) Ausgaben sheet: insert two new expense rows (Filmcoopi Zuerich AG /
#    "Rechnung Ingeborg Bachmann" invoice) just above the existing
#    "Personalaufwand" salary-payment block, and grow Table16 to include
#    them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Einnahmen: new row at the bottom of Table3 (A1:G7 -> A1:G8)
# ---------------------------------------------------------------------
$wsIn = $wb.Worksheets.Item("Einnahmen")
$loIn = $wsIn.ListObjects.Item(1)
$newIncomeRow = $loIn.ListRows.Add()
$rowRangeIn = $newIncomeRow.Range

$rowRangeIn.Cells.Item(1, 1).Value = "Sonstige Einnahmen"
$rowRangeIn.Cells.Item(1, 2).Value = "Aargauer Kuratorium"
$rowRangeIn.Cells.Item(1, 3).Value = 45265
$rowRangeIn.Cells.Item(1, 4).Value = 12000
$rowRangeIn.Cells.Item(1, 5).Value = "Aargauer Kuratorium"

# ---------------------------------------------------------------------
# Ausgaben: drop the old row 16 (the 25.3.2023 Personalaufwand entry is
# superseded / no longer listed), then insert 2 fresh rows above what is
# now the first remaining Personalaufwand / Gehaltszahlung Projektleitung
# entry, fill them in with the new Filmcoopi invoice, and grow Table16 so
# the inserted rows become part of the table.
# ---------------------------------------------------------------------
$wsOut = $wb.Worksheets.Item("Ausgaben")
$loOut = $wsOut.ListObjects.Item(1)

$wsOut.Rows.Item(16).Delete()

$wsOut.Rows.Item(16).Resize(2).Insert()

$wsOut.Cells.Item(16, 1).Value = "Film"
$wsOut.Cells.Item(16, 2).Value = "Rechnung Ingeborg Bachmann"
$wsOut.Cells.Item(16, 3).Value = 45291
$wsOut.Cells.Item(16, 4).Value = 127.05
$wsOut.Cells.Item(16, 5).Value = "Filmcoopi Zürich AG"
$wsOut.Cells.Item(16, 6).Value = "Heinrichstrasse 114, 8005 Zürich"
$wsOut.Cells.Item(16, 7).Value = "93 64950 00000 00000 00008 38954"
$wsOut.Cells.Item(16, 8).Value = "00083895"
$wsOut.Cells.Item(16, 9).Value = 45267

$wsOut.Cells.Item(17, 1).Value = "Film"
$wsOut.Cells.Item(17, 2).Value = "Rechnung Ingeborg Bachmann"
$wsOut.Cells.Item(17, 3).Value = 45291
$wsOut.Cells.Item(17, 4).Value = 231.55
$wsOut.Cells.Item(17, 5).Value = "Filmcoopi Zürich AG"
$wsOut.Cells.Item(17, 6).Value = "Heinrichstrasse 114, 8005 Zürich"
$wsOut.Cells.Item(17, 7).Value = "93 64950 00000 00000 00008 38946"
$wsOut.Cells.Item(17, 8).Value = "00083894"
$wsOut.Cells.Item(17, 9).Value = 45256

$loOut.Resize($wsOut.Range("A1:I26"))
